$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "26.827.87"
$r.Style = "Normal"
$ws.Range("E2").Value = "  -0.79%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.813.31"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +0.82%  "
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.005"
$r.Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "309.10"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  -0.07%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.4325"
$r.Style = "Normal"
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.3709"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +3.14%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.07257"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +2.77%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "20.90"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +3.21%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "1.943.90"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +3.86%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "6.643"
$r.Style = "Normal"
$ws.Range("E13").Value = "  +4.39%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "5.359"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +1.40%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.06925"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "1.005"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -0.43%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "80.62"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.02%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.000008876"
$r.Style = "Normal"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +1.23%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "26.859.20"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -1.60%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "5.217"
$r.Style = "Normal"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("E23").Value = "  +0.81%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "2.172.19"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +4.23%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "154.04"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "1.872"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -4.37%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "18.29"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "5.232"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +4.35%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "1.902"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +15.06%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "115.38"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.08945"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -0.48%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "0.7566"
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.174"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +7.40%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "4.437"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +1.94%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "2.810"
$r.Style = "Normal"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  +0.25%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "1.132"
$r.Style = "Normal"
$ws.Range("E37").Value = "  +4.97%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.05234"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +1.65%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.01928"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +1.22%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.5087"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +2.25%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.1649"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +1.06%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "2.692"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +2.84%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "6.555"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +9.78%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "8.301"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +2.66%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "10.48"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +1.90%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "106.81"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  -0.06%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "1.657"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +3.42%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.4584"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.06286"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.54%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "1.814"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +5.08%  "
